$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width adjustments (swap L and P widths)
$ws.Range("L1").ColumnWidth = 10.75
$ws.Range("P1").ColumnWidth = 11.75

# Cell value updates
$ws.Range("AN1").Value = 0.81853653715329111
$ws.Range("AW1").Value = 0.98415459871931021
$ws.Range("I2").Value = 0.96356190245638607
$ws.Range("L2").Value = 0.79606282101327708
$ws.Range("AN3").Value = 0.67882703917058218
$ws.Range("BH3").Value = 0.83153202074567711
$ws.Range("BL3").Value = 0.95145153373640112
$ws.Range("BH4").Value = 0.70598460776351513
$ws.Range("D5").Value = 0.99968334510845669
$ws.Range("G5").Value = 0.89426261388626993
$ws.Range("AW5").Value = 0.85751345853561056
$ws.Range("BL5").Value = 0.89651266952284736
$ws.Range("N6").Value = 0.69023304643924432
$ws.Range("AB6").Value = 0.6924828255848603
$ws.Range("BD6").Value = 0.61567563810087267
$ws.Range("C7").Value = 0.60396962950220157
$ws.Range("I7").Value = 0.73359504870179193
$ws.Range("AO7").Value = 0.98885093796562251
$ws.Range("AV7").Value = 0.93821568609995976
$ws.Range("BF7").Value = 0.7521630347596826
$ws.Range("W8").Value = 0.93189145405603013
$ws.Range("AG8").Value = 0.95375649082006231
$ws.Range("AO9").Value = 0.86626409877959176
$ws.Range("U10").Value = 0.65684634720302748
$ws.Range("AF10").Value = 0.88631403877037973
$ws.Range("I11").Value = 0.9657028963804084
$ws.Range("BO11").Value = 0.94290284940867675
$ws.Range("Q12").Value = 0.99483331603651237
$ws.Range("Q14").Value = 0.80513371023499691
$ws.Range("AW14").Value = 0.98574159410724804
$ws.Range("H15").Value = 0.87914054540852749
$ws.Range("M15").Value = 0.87921668665707187
$ws.Range("K16").Value = 0.85504902102714619
$ws.Range("P18").Value = 0.69495220866761243
$ws.Range("Q18").Value = 0.96187811998102468
$ws.Range("BG18").Value = 0.78298824860910299
$ws.Range("R19").Value = 0.83380268156851889
$ws.Range("R20").Value = 0.81550733653107987
$ws.Range("X20").Value = 0.86499840538991823
$ws.Range("W21").Value = 0.90483834457349355
$ws.Range("AX21").Value = 0.92525713286973565
$ws.Range("S22").Value = 0.83729961157887434
$ws.Range("T22").Value = 0.93787341877600627
$ws.Range("U22").Value = 0.99679695095948717
$ws.Range("W22").Value = 0.93399920526337121
$ws.Range("AX22").Value = 0.77518701855719896
$ws.Range("AU23").Value = 0.97778533225775632
$ws.Range("W24").Value = 0.96770355764718841
$ws.Range("Y24").Value = 0.6651546977563394
$ws.Range("AO24").Value = 0.99324639509202073
$ws.Range("F25").Value = 0.85626572366756415
$ws.Range("AB25").Value = 0.8829747584210399
$ws.Range("BN25").Value = 0.87064683267847298
$ws.Range("A26").Value = 0.9021129836875792
$ws.Range("BH26").Value = 0.9184152565874566
$ws.Range("B27").Value = 0.97102022645486441
$ws.Range("AB27").Value = 0.87503720959953335
$ws.Range("S28").Value = 0.99675217194854171
$ws.Range("AZ28").Value = 0.60660899351899611
$ws.Range("K29").Value = 0.88052299056447336
$ws.Range("BH29").Value = 0.98033916927063425
$ws.Range("B30").Value = 0.98250211512961483
$ws.Range("I30").Value = 0.69819336875736504
$ws.Range("X30").Value = 0.8407202983599924
$ws.Range("T31").Value = 0.95180050287384343
$ws.Range("F32").Value = 0.76584547244326884
$ws.Range("AE32").Value = 0.90600136591749347
$ws.Range("BH32").Value = 0.97807449514275291
$ws.Range("M33").Value = 0.5966393349801089
$ws.Range("R33").Value = 0.59532634769114967
$ws.Range("BH33").Value = 0.65774047806490887
$ws.Range("O34").Value = 0.98732828631048919
$ws.Range("Q34").Value = 0.72391854558777946
$ws.Range("AU34").Value = 0.69337906408763184
$ws.Range("BG34").Value = 0.7585611763813398
$ws.Range("AH35").Value = 0.88464938709516616
$ws.Range("AK35").Value = 0.72228005828806463
$ws.Range("G36").Value = 0.99190243936976752
$ws.Range("BP36").Value = 0.94033240228524861
$ws.Range("AU37").Value = 0.9750061174573128
$ws.Range("BP37").Value = 0.89998868464203508
$ws.Range("B38").Value = 0.83481665652698922
$ws.Range("AJ38").Value = 0.76628318684451879
$ws.Range("BF38").Value = 0.75160396495374415
$ws.Range("AQ39").Value = 0.91226631968519922
$ws.Range("J41").Value = 0.67441895614693093
$ws.Range("AD41").Value = 0.78509165325295704
$ws.Range("AB42").Value = 0.75184839117753111
$ws.Range("AT42").Value = 0.95765584185953267
$ws.Range("AZ42").Value = 0.706846236567783
$ws.Range("J43").Value = 0.95824014750046704
$ws.Range("AR43").Value = 0.83309068545579268
$ws.Range("AC44").Value = 0.79458756551137411
$ws.Range("AD44").Value = 0.96273902513736087
$ws.Range("AP44").Value = 0.64530449860681649
$ws.Range("W45").Value = 0.68248932884269187
$ws.Range("BA45").Value = 0.69192825369934119
$ws.Range("BB45").Value = 0.74468995047176012
$ws.Range("AZ46").Value = 0.87963252685197635
$ws.Range("C48").Value = 0.85422382291063004
$ws.Range("BA48").Value = 0.92389137793984932
$ws.Range("AU49").Value = 0.92182842838517631
$ws.Range("BJ49").Value = 0.9180578086658604
$ws.Range("AZ51").Value = 0.88917767511714718
$ws.Range("BF51").Value = 0.57969283740644395
$ws.Range("BM52").Value = 0.69785637809723577
$ws.Range("AM53").Value = 0.9299240804615404
$ws.Range("BC53").Value = 0.89737708558275597
$ws.Range("BN54").Value = 0.94605729318944165
$ws.Range("U55").Value = 0.88711980241992827
$ws.Range("AI55").Value = 0.82631556512283399
$ws.Range("Q56").Value = 0.88750933566758461
$ws.Range("Z56").Value = 0.90206090561151264
$ws.Range("G57").Value = 0.69203475996501362
$ws.Range("AL57").Value = 0.67314660260907166
$ws.Range("AD58").Value = 0.98131919256801758
$ws.Range("AP59").Value = 0.89938087861165228
$ws.Range("T60").Value = 0.83503748758083152
$ws.Range("U61").Value = 0.96533577026540307
$ws.Range("AN61").Value = 0.82861537907834215
$ws.Range("BB61").Value = 0.86656210227339869
$ws.Range("BI62").Value = 0.92973249849294948
$ws.Range("BK62").Value = 0.77075295597894078
$ws.Range("B63").Value = 0.66940990591836269
$ws.Range("O63").Value = 0.97163835790367747
$ws.Range("AK63").Value = 0.91810027656034721
$ws.Range("AU63").Value = 0.97547184182106683
$ws.Range("AT64").Value = 0.90145439167467722
$ws.Range("K65").Value = 0.61471532761856673
$ws.Range("V65").Value = 0.8361581678641119
$ws.Range("BA65").Value = 0.882864324082
$ws.Range("BE66").Value = 0.99850804854251396
$ws.Range("AT67").Value = 0.82747518917787533
$ws.Range("BE67").Value = 0.77206603345776492
$ws.Range("E68").Value = 0.9464508020095711
$ws.Range("AL68").Value = 0.58042667229973599
